$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2556606666666667
$ws.Range("H2").Value = 0.7669820000000001
$ws.Range("I2").Value = 0.01354513404628681
$ws.Range("J2").Value = 0.01354513404628681
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 93.78975666666668
$ws.Range("N2").Value = 281.36927
$ws.Range("O2").Value = 0.997863063099077
$ws.Range("P2").Value = 0.9978630630990771
$ws.Range("Q2").Value = 23.97835171590445
$ws.Range("R2").Value = 215.80516544314
$ws.Range("S2").Value = 0.01351618894951535
$ws.Range("T2").Value = 0.01351618894951535

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2556606666666667
$ws.Range("H3").Value = 0.7669820000000001
$ws.Range("I3").Value = 0.01354513404628681
$ws.Range("J3").Value = 0.01354513404628681
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.09680433333333333
$ws.Range("N3").Value = 0.290413
$ws.Range("O3").Value = 0.001029936231997873
$ws.Range("P3").Value = 0.001029936231997873
$ws.Range("Q3").Value = 0.02474906039622222
$ws.Range("R3").Value = 0.222741543566
$ws.Range("S3").Value = 0.00001395062432153874
$ws.Range("T3").Value = 0.00001395062432153874

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2556606666666667
$ws.Range("H4").Value = 0.7669820000000001
$ws.Range("I4").Value = 0.01354513404628681
$ws.Range("J4").Value = 0.01354513404628681
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.05343666666666667
$ws.Range("N4").Value = 0.16031
$ws.Range("O4").Value = 0.0005685319780849309
$ws.Range("P4").Value = 0.000568531978084931
$ws.Range("Q4").Value = 0.01366165382444445
$ws.Range("R4").Value = 0.12295488442
$ws.Range("S4").Value = 0.000007700841852760982
$ws.Range("T4").Value = 0.000007700841852760984

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2556606666666667
$ws.Range("H5").Value = 0.7669820000000001
$ws.Range("I5").Value = 0.01354513404628681
$ws.Range("J5").Value = 0.01354513404628681
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.050611
$ws.Range("N5").Value = 0.151833
$ws.Range("O5").Value = 0.0005384686908400556
$ws.Range("P5").Value = 0.0005384686908400557
$ws.Range("Q5").Value = 0.01293924200066667
$ws.Range("R5").Value = 0.116453178006
$ws.Range("S5").Value = 0.000007293630597157121
$ws.Range("T5").Value = 0.000007293630597157123

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 16.677026
$ws.Range("H6").Value = 50.03107799999999
$ws.Range("I6").Value = 0.8835639662863414
$ws.Range("J6").Value = 0.8835639662863415
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 93.78975666666668
$ws.Range("N6").Value = 281.36927
$ws.Range("O6").Value = 0.997863063099077
$ws.Range("P6").Value = 0.9978630630990771
$ws.Range("Q6").Value = 1564.134210463673
$ws.Range("R6").Value = 14077.20789417306
$ws.Range("S6").Value = 0.8816758458424582
$ws.Range("T6").Value = 0.8816758458424584

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 16.677026
$ws.Range("H7").Value = 50.03107799999999
$ws.Range("I7").Value = 0.8835639662863414
$ws.Range("J7").Value = 0.8835639662863415
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.09680433333333333
$ws.Range("N7").Value = 0.290413
$ws.Range("O7").Value = 0.001029936231997873
$ws.Range("P7").Value = 0.001029936231997873
$ws.Range("Q7").Value = 1.614408383912666
$ws.Range("R7").Value = 14.529675455214
$ws.Range("S7").Value = 0.0009100145421660502
$ws.Range("T7").Value = 0.0009100145421660506

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 16.677026
$ws.Range("H8").Value = 50.03107799999999
$ws.Range("I8").Value = 0.8835639662863414
$ws.Range("J8").Value = 0.8835639662863415
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.05343666666666667
$ws.Range("N8").Value = 0.16031
$ws.Range("O8").Value = 0.0005685319780849309
$ws.Range("P8").Value = 0.000568531978084931
$ws.Range("Q8").Value = 0.8911646793533332
$ws.Range("R8").Value = 8.02048211418
$ws.Range("S8").Value = 0.0005023343695173409
$ws.Range("T8").Value = 0.000502334369517341

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 16.677026
$ws.Range("H9").Value = 50.03107799999999
$ws.Range("I9").Value = 0.8835639662863414
$ws.Range("J9").Value = 0.8835639662863415
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.050611
$ws.Range("N9").Value = 0.151833
$ws.Range("O9").Value = 0.0005384686908400556
$ws.Range("P9").Value = 0.0005384686908400557
$ws.Range("Q9").Value = 0.8440409628859998
$ws.Range("R9").Value = 7.596368665973999
$ws.Range("S9").Value = 0.0004757715321996533
$ws.Range("T9").Value = 0.0004757715321996534

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.942037333333333
$ws.Range("H10").Value = 5.826112
$ws.Range("I10").Value = 0.1028908996673717
$ws.Range("J10").Value = 0.1028908996673717
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 93.78975666666668
$ws.Range("N10").Value = 281.36927
$ws.Range("O10").Value = 0.997863063099077
$ws.Range("P10").Value = 0.9978630630990771
$ws.Range("Q10").Value = 182.1432089309156
$ws.Range("R10").Value = 1639.28888037824
$ws.Range("S10").Value = 0.1026710283071034
$ws.Range("T10").Value = 0.1026710283071034

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.942037333333333
$ws.Range("H11").Value = 5.826112
$ws.Range("I11").Value = 0.1028908996673717
$ws.Range("J11").Value = 0.1028908996673717
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.09680433333333333
$ws.Range("N11").Value = 0.290413
$ws.Range("O11").Value = 0.001029936231997873
$ws.Range("P11").Value = 0.001029936231997873
$ws.Range("Q11").Value = 0.1879976293617778
$ws.Range("R11").Value = 1.691978664256
$ws.Range("S11").Value = 0.0001059710655102841
$ws.Range("T11").Value = 0.0001059710655102841

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.942037333333333
$ws.Range("H12").Value = 5.826112
$ws.Range("I12").Value = 0.1028908996673717
$ws.Range("J12").Value = 0.1028908996673717
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.05343666666666667
$ws.Range("N12").Value = 0.16031
$ws.Range("O12").Value = 0.0005685319780849309
$ws.Range("P12").Value = 0.000568531978084931
$ws.Range("Q12").Value = 0.1037760016355556
$ws.Range("R12").Value = 0.93398401472
$ws.Range("S12").Value = 0.00005849676671482901
$ws.Range("T12").Value = 0.00005849676671482903

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.942037333333333
$ws.Range("H13").Value = 5.826112
$ws.Range("I13").Value = 0.1028908996673717
$ws.Range("J13").Value = 0.1028908996673717
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.050611
$ws.Range("N13").Value = 0.151833
$ws.Range("O13").Value = 0.0005384686908400556
$ws.Range("P13").Value = 0.0005384686908400557
$ws.Range("Q13").Value = 0.09828845147733332
$ws.Range("R13").Value = 0.884596063296
$ws.Range("S13").Value = 0.00005540352804324517
$ws.Range("T13").Value = 0.00005540352804324519
